$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.882.49'
$ws.Range("E2").Value = '  +7.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.686.33'
$ws.Range("E3").Value = '  +11.48%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '511.29'
$ws.Range("E5").Value = '  +4.70%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.02'
$ws.Range("E6").Value = '  +2.13%  '

# Row 7
$ws.Range("E7").Value = '  +0.36%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.606'
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.681.36'
$ws.Range("E9").Value = '  +10.50%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.47'
$ws.Range("E10").Value = '  +3.13%  '

# Row 11
$ws.Range("E11").Value = '  +4.90%  '

# Row 12
$ws.Range("E12").Value = '  +3.79%  '

# Row 13
$ws.Range("E13").Value = '  +0.97%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.148.76'
$ws.Range("E14").Value = '  +11.25%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.936.00'
$ws.Range("E15").Value = '  +7.13%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.78'
$ws.Range("E16").Value = '  +4.69%  '

# Row 17
$ws.Range("E17").Value = '  +4.89%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.682.14'
$ws.Range("E18").Value = '  +10.68%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.80'
$ws.Range("E19").Value = '  +1.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.80'
$ws.Range("E20").Value = '  +7.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.53'
$ws.Range("E21").Value = '  +5.29%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.20'
$ws.Range("E22").Value = '  +2.95%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.45'
$ws.Range("E24").Value = '  +3.85%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  +3.84%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.799.91'
$ws.Range("E26").Value = '  +11.29%  '

# Row 27
$ws.Range("E27").Value = '  +3.35%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.29%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0860'
$ws.Range("E29").Value = '  +9.66%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.53'
$ws.Range("E30").Value = '  +2.70%  '

# Row 31
$ws.Range("E31").Value = '  +0.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.77'
$ws.Range("E32").Value = '  +4.95%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.62'
$ws.Range("E33").Value = '  +5.73%  '

# Row 34
$ws.Range("E34").Value = '  +3.49%  '

# Row 35
$ws.Range("E35").Value = '  +7.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.08'
$ws.Range("E36").Value = '  +9.48%  '

# Row 37
$ws.Range("E37").Value = '  +5.00%  '

# Row 38
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '311.64'
$ws.Range("E38").Value = '  +16.28%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  +10.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.856'
$ws.Range("E40").Value = '  +1.18%  '

# Row 41
$ws.Range("E41").Value = '  +6.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.837'
$ws.Range("E42").Value = '  +30.26%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.46'
$ws.Range("E43").Value = '  +3.74%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.645'
$ws.Range("E44").Value = '  +8.71%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0576'
$ws.Range("E45").Value = '  +7.81%  '

# Row 46
$ws.Range("E46").Value = '  -0.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.58%  '

# Row 48
$ws.Range("E48").Value = '  +14.33%  '

# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.063.50'
$ws.Range("E49").Value = '  +10.16%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0237'
$ws.Range("E50").Value = '  +3.46%  '

# Row 51
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.85'
$ws.Range("E51").Value = '  +4.38%  '

Write-Host "Done"